$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 67.8679804978435
$ws.Range("C2").Value = 66.47906002391
$ws.Range("D2").Value = 69.256900971777
$ws.Range("C3").Value = 63.3414374246431
$ws.Range("D3").Value = 68.0972229127337
$ws.Range("B11").Value = 67.987520077397
$ws.Range("C11").Value = 64.4000884639409
$ws.Range("D11").Value = 71.5749516908532
$ws.Range("B12").Value = 74.0846263450239
$ws.Range("C12").Value = 70.2362618026456
$ws.Range("D12").Value = 77.9329908874022
$ws.Range("B14").Value = 64.0735542026883
$ws.Range("C14").Value = 61.0923553714412
$ws.Range("D14").Value = 67.0547530339354
$ws.Range("B15").Value = 55.3718303775116
$ws.Range("C15").Value = 54.4833021444574
$ws.Range("D15").Value = 56.2603586105659
$ws.Range("C16").Value = 61.1949540788548
$ws.Range("D16").Value = 64.02596086048
$ws.Range("C17").Value = 58.8593519707779
$ws.Range("D17").Value = 65.5953518572809
$ws.Range("B18").Value = 56.7933697165468
$ws.Range("C18").Value = 53.7437040271019
$ws.Range("D18").Value = 59.8430354059916
$ws.Range("C19").Value = 72.6197742299908
$ws.Range("D19").Value = 78.8410281401418
$ws.Range("C20").Value = 48.5071840403106
$ws.Range("D20").Value = 57.6094230654795
$ws.Range("C21").Value = 65.1502269043595
$ws.Range("D21").Value = 72.6616418035527
$ws.Range("C22").Value = 55.8991604374395
$ws.Range("D22").Value = 62.73999632971
$ws.Range("C23").Value = 53.1374134879005
$ws.Range("D23").Value = 71.5777265405306
$ws.Range("B24").Value = 54.8090076208069
$ws.Range("C24").Value = 52.3361474530153
$ws.Range("D24").Value = 57.2818677885985
$ws.Range("B25").Value = 49.0282605200155
$ws.Range("C25").Value = 46.6386706813937
$ws.Range("D25").Value = 51.4178503586373
$ws.Range("C26").Value = 80.8509054008945
$ws.Range("D26").Value = 88.4348158368641
$ws.Range("B27").Value = 47.4912333237318
$ws.Range("C27").Value = 45.9334227038904
$ws.Range("D27").Value = 49.0490439435732
$ws.Range("B28").Value = 64.7240480131028
$ws.Range("C28").Value = 63.0045795926267
$ws.Range("D28").Value = 66.4435164335789
$ws.Range("C29").Value = 66.154032049576
$ws.Range("D29").Value = 71.4780870671454
$ws.Range("B37").Value = 66.5933688883394
$ws.Range("C37").Value = 62.4403846149291
$ws.Range("D37").Value = 70.7463531617498
$ws.Range("B38").Value = 57.2690723237937
$ws.Range("C38").Value = 50.7470363646041
$ws.Range("D38").Value = 63.7911082829833
$ws.Range("B40").Value = 60.3471262597791
$ws.Range("C40").Value = 56.8080642142266
$ws.Range("D40").Value = 63.8861883053316
$ws.Range("B41").Value = 67.2054760884641
$ws.Range("C41").Value = 65.5001731723724
$ws.Range("D41").Value = 68.9107790045557
$ws.Range("B42").Value = 73.150166637206
$ws.Range("C42").Value = 70.5221212688109
$ws.Range("D42").Value = 75.7782120056012
$ws.Range("B50").Value = 67.292535125419
$ws.Range("C50").Value = 63.1364588496784
$ws.Range("D50").Value = 71.4486114011597
$ws.Range("B51").Value = 63.9033952437169
$ws.Range("C51").Value = 58.2061893987969
$ws.Range("D51").Value = 69.6006010886369
$ws.Range("B53").Value = 58.9468339412676
$ws.Range("C53").Value = 55.254083851168
$ws.Range("D53").Value = 62.6395840313673
$ws.Range("B54").Value = 68.518253119246
$ws.Range("C54").Value = 66.8882724435162
$ws.Range("D54").Value = 70.1482337949759
$ws.Range("C55").Value = 70.7689055039649
$ws.Range("D55").Value = 75.9022703411351
$ws.Range("B63").Value = 70.8149579691965
$ws.Range("C63").Value = 66.8120913855674
$ws.Range("D63").Value = 74.8178245528256
$ws.Range("B64").Value = 60.4866837966034
$ws.Range("C64").Value = 54.4980227132752
$ws.Range("D64").Value = 66.4753448799316
$ws.Range("B66").Value = 63.6828257869627
$ws.Range("C66").Value = 60.3796020159564
$ws.Range("D66").Value = 66.986049557969
$ws.Range("B67").Value = 61.895417279981
$ws.Range("C67").Value = 60.1645023547081
$ws.Range("D67").Value = 63.6263322052539
$ws.Range("C68").Value = 65.1898278582191
$ws.Range("D68").Value = 70.736821035104
$ws.Range("B76").Value = 63.3861854984241
$ws.Range("C76").Value = 59.2235413127094
$ws.Range("D76").Value = 67.5488296841387
$ws.Range("B77").Value = 54.11746546684
$ws.Range("C77").Value = 48.0133810262521
$ws.Range("D77").Value = 60.2215499074278
$ws.Range("B79").Value = 56.0767240057917
$ws.Range("C79").Value = 52.589586582689
$ws.Range("D79").Value = 59.5638614288945
